# Apply scheduled-runner updates to Masamune_Profits workbook (per-sheet leve profit recalcs)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 14153.411
$ws.Range("I132").Value = 2023.375
$ws.Range("J132").Value = 100411.445
$ws.Range("K132").Value = 6070.125
$ws.Range("L132").Value = 301234.335
$ws.Range("M132").Value = -3540.125
$ws.Range("N132").Value = -306294.335

$ws.Range("H133").Value = 47531.43
$ws.Range("J133").Value = 47531.43
$ws.Range("L133").Value = 47531.43
$ws.Range("N133").Value = -57651.43

$ws.Range("H134").Value = 30863.158
$ws.Range("J134").Value = 30863.158
$ws.Range("L134").Value = 30863.158
$ws.Range("N134").Value = -41003.158

$ws.Range("H136").Value = 37400
$ws.Range("J136").Value = 37400
$ws.Range("L136").Value = 37400
$ws.Range("N136").Value = -47600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36412.164
$ws.Range("I32").Value = 43763.23
$ws.Range("J32").Value = 18034.5
$ws.Range("K32").Value = 43763.23
$ws.Range("L32").Value = 18034.5
$ws.Range("M32").Value = -43476.23
$ws.Range("N32").Value = -18608.5

$ws.Range("H104").Value = 35041.75
$ws.Range("J104").Value = 35041.75
$ws.Range("L104").Value = 35041.75
$ws.Range("N104").Value = -42029.75

$ws.Range("H122").Value = 2230.7585
$ws.Range("I122").Value = 2474.6
$ws.Range("J122").Value = 1688.8889
$ws.Range("K122").Value = 7423.799999999999
$ws.Range("L122").Value = 5066.6667
$ws.Range("M122").Value = -4973.799999999999
$ws.Range("N122").Value = -9966.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 27975
$ws.Range("J106").Value = 27975
$ws.Range("L106").Value = 27975
$ws.Range("N106").Value = -30499

$ws.Range("H111").Value = 22680
$ws.Range("J111").Value = 22680
$ws.Range("L111").Value = 22680
$ws.Range("N111").Value = -30860

$ws.Range("H115").Value = 34249.332
$ws.Range("J115").Value = 34249.332
$ws.Range("L115").Value = 34249.332
$ws.Range("N115").Value = -36599.332

$ws.Range("H116").Value = 37549.75
$ws.Range("J116").Value = 37549.75
$ws.Range("L116").Value = 37549.75
$ws.Range("N116").Value = -46727.75

$ws.Range("H118").Value = 33969.25
$ws.Range("J118").Value = 33969.25
$ws.Range("L118").Value = 33969.25
$ws.Range("N118").Value = -37283.25

$ws.Range("H119").Value = 48504.332
$ws.Range("J119").Value = 48504.332
$ws.Range("L119").Value = 48504.332
$ws.Range("N119").Value = -58180.332

$ws.Range("H122").Value = 63902.58
$ws.Range("I122").Value = 93107.234
$ws.Range("J122").Value = 625.8333
$ws.Range("K122").Value = 279321.702
$ws.Range("L122").Value = 1877.4999
$ws.Range("M122").Value = -276871.702
$ws.Range("N122").Value = -6777.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 125001896
$ws.Range("J9").Value = 2678
$ws.Range("L9").Value = 8034
$ws.Range("N9").Value = -8482

$ws.Range("H22").Value = 21039
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 21039
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 63117
$ws.Range("N22").Value = -63455
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 21039
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 21039
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 63117
$ws.Range("N27").Value = -63321
$ws.Range("M27").ClearContents()

$ws.Range("H76").Value = 4659.553
$ws.Range("I76").Value = 1950
$ws.Range("J76").Value = 4779.9775
$ws.Range("K76").Value = 5850
$ws.Range("L76").Value = 14339.9325
$ws.Range("M76").Value = -5467
$ws.Range("N76").Value = -15105.9325

$ws.Range("H79").Value = 4659.553
$ws.Range("I79").Value = 1950
$ws.Range("J79").Value = 4779.9775
$ws.Range("K79").Value = 5850
$ws.Range("L79").Value = 14339.9325
$ws.Range("M79").Value = -4524
$ws.Range("N79").Value = -16991.9325

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4374.7
$ws.Range("I107").Value = 749.4
$ws.Range("J107").Value = 8000
$ws.Range("K107").Value = 749.4
$ws.Range("L107").Value = 8000
$ws.Range("M107").Value = 1170.6
$ws.Range("N107").Value = -11840

$ws.Range("H130").Value = 45303.7
$ws.Range("J130").Value = 45303.7
$ws.Range("L130").Value = 45303.7
$ws.Range("N130").Value = -55343.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 30390.125
$ws.Range("J105").Value = 30390.125
$ws.Range("L105").Value = 30390.125
$ws.Range("N105").Value = -37378.125

$ws.Range("H106").Value = 29287.334
$ws.Range("J106").Value = 29287.334
$ws.Range("L106").Value = 29287.334
$ws.Range("N106").Value = -31811.334

$ws.Range("H110").Value = 34808
$ws.Range("J110").Value = 34808
$ws.Range("L110").Value = 34808
$ws.Range("N110").Value = -42988

$ws.Range("H114").Value = 19158
$ws.Range("J114").Value = 19158
$ws.Range("L114").Value = 19158
$ws.Range("N114").Value = -27836

$ws.Range("H116").Value = 31668
$ws.Range("J116").Value = 31668
$ws.Range("L116").Value = 31668
$ws.Range("N116").Value = -40846

$ws.Range("H121").Value = 16649.5
$ws.Range("J121").Value = 16649.5
$ws.Range("L121").Value = 16649.5
$ws.Range("N121").Value = -20143.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 35308
$ws.Range("J16").Value = 35308
$ws.Range("L16").Value = 35308
$ws.Range("N16").Value = -35892

$ws.Range("H106").Value = 31801.777
$ws.Range("J106").Value = 31801.777
$ws.Range("L106").Value = 31801.777
$ws.Range("N106").Value = -34325.777

$ws.Range("H110").Value = 21637.5
$ws.Range("J110").Value = 21637.5
$ws.Range("L110").Value = 21637.5
$ws.Range("N110").Value = -29817.5

$ws.Range("H114").Value = 23392
$ws.Range("J114").Value = 23392
$ws.Range("L114").Value = 23392
$ws.Range("N114").Value = -32070

$ws.Range("H121").Value = 29231.25
$ws.Range("J121").Value = 29231.25
$ws.Range("L121").Value = 29231.25
$ws.Range("N121").Value = -32725.25

$ws.Range("H132").Value = 2066.0857
$ws.Range("I132").Value = 1510.75
$ws.Range("J132").Value = 2806.5334
$ws.Range("K132").Value = 4532.25
$ws.Range("L132").Value = 8419.600199999999
$ws.Range("M132").Value = -2002.25
$ws.Range("N132").Value = -13479.6002
